$wb = $excel.ActiveWorkbook

# --- Sheet 1: StrataSystem ---
$ws1 = $wb.Worksheets.Item("StrataSystem")

# Fix typo: Stine Frie -> Stine Fire (rows 6 and 7, PersonInCharge column)
$ws1.Range("F6").Value = "Stine Fire"
$ws1.Range("F7").Value = "Stine Fire"

# Add new StrataSystem row 8: WGINOR
$ws1.Range("A8").Value = "WGINOR"
$ws1.Range("B8").Value = "WGINOR"
$ws1.Range("C8").Value = 2023
$ws1.Range("D8").Value = "Suggested polygons"
$ws1.Range("E8").Value = "NA"
$ws1.Range("F8").Value = "Solva Eliassen"
$ws1.Range("H8").Value = "NA"

# --- Sheet 2: Strata ---
$ws2 = $wb.Worksheets.Item("Strata")

# Add new Strata rows 93-102 for the WGINOR strata system
$ws2.Range("A93").Value = "WGINORAW"
$ws2.Range("B93").Value = "WGINOR"
$ws2.Range("C93").Value = 'AW'
$ws2.Range("A94").Value = "WGINOREIW"
$ws2.Range("B94").Value = "WGINOR"
$ws2.Range("C94").Value = 'EIW'
$ws2.Range("A95").Value = "WGINORFaroe"
$ws2.Range("B95").Value = "WGINOR"
$ws2.Range("C95").Value = 'Faroe'
$ws2.Range("A96").Value = "WGINORGS"
$ws2.Range("B96").Value = "WGINOR"
$ws2.Range("C96").Value = 'GS'
$ws2.Range("A97").Value = "WGINORIFF"
$ws2.Range("B97").Value = "WGINOR"
$ws2.Range("C97").Value = 'IFF'
$ws2.Range("A98").Value = "WGINORLBE"
$ws2.Range("B98").Value = "WGINOR"
$ws2.Range("C98").Value = 'LBE'
$ws2.Range("A99").Value = "WGINORLBW"
$ws2.Range("B99").Value = "WGINOR"
$ws2.Range("C99").Value = 'LBW'
$ws2.Range("A100").Value = "WGINORNorSeaE"
$ws2.Range("B100").Value = "WGINOR"
$ws2.Range("C100").Value = 'NorSeaE'
$ws2.Range("A101").Value = "WGINORNorSeaW"
$ws2.Range("B101").Value = "WGINOR"
$ws2.Range("C101").Value = 'NorSeaW'
$ws2.Range("A102").Value = "WGINORRH"
$ws2.Range("B102").Value = "WGINOR"
$ws2.Range("C102").Value = 'RH'

# --- Sheet 4: Coordinates ---
$ws4 = $wb.Worksheets.Item("Coordinates")

# Add new Coordinates rows 93-102 with the WGINOR polygon WKT data
$ws4.Range("A93").Value = "WGINORAW"
$ws4.Range("B93").Value = 'POLYGON((-25.7999992370605  60.5,-24.8999996185303  60.7333335876465,-24.2333335876465  60.7666664123535,-23.7999992370605  60.7000007629395,-23.5  60.533332824707,-22.033332824707  61.1333351135254,-21.3333339691162  61.2000007629395,-21.1333332061768  61.3666648864746,-21.2666664123535  61.533332824707,-20.7000007629395  61.7666664123535,-20.2000007629395  61.6666679382324,-19.6000003814697  61.7333335876465,-18.966667175293  61.6666679382324,-18.966667175293  61.966667175293,-18.6666660308838  61.966667175293,-18.0666675567627  61.8666648864746,-17.6333332061768  62.0666656494141,-17.533332824707  62.3333320617676,-17  62.5999984741211,-16.6333332061768  62.5,-15.5666666030884  62.6333351135254,-15.2666664123535  62.7333335876465,-14.7666664123535  61.7666664123535,-14.3999996185303  61.5,-14.6000003814697  61.2999992370605,-15  60.2333335876465,-15.6666669845581  59.966667175293,-16.0666675567627  59.7999992370605,-16.5666675567627  59.8333320617676,-18.1333332061768  59.4000015258789,-19.466667175293  58.7999992370605,-19.7999992370605  58.533332824707,-20.7333335876465  57.1666679382324,-21.033332824707  57.0666656494141,-22.1000003814697  56.466667175293,-22.5666675567627  56.4000015258789,-23.1666660308838  55.9333343505859,-23.2999992370605  55.7999992370605,-22.466667175293  55.4333343505859,-21.8666667938232  55.4333343505859,-20.8666667938232  55.5999984741211,-20.1333332061768  55.1666679382324,-20.2999992370605  55.033332824707,-20.1666660308838  53.7666664123535,-19.8666667938232  53.7000007629395,-18.6666660308838  54.1333351135254,-18.1333332061768  54.2999992370605,-17.6666660308838  54.5,-17.3333339691162  54.5666656494141,-16.533332824707  55.0666656494141,-15.7333335876465  55.2333335876465,-15.0333337783813  55.533332824707,-14.1666669845581  56,-12.8666667938232  57,-11.7333335876465  57.2333335876465,-11.5333337783813  57.5999984741211,-11.1666669845581  58,-10.6666669845581  58.2333335876465,-10.0666666030884  58.0999984741211,-9.53333377838135  58.2000007629395,-8.80000019073486  58.2999992370605,-8.33333301544189  58.3666648864746,-7.69999980926514  58.6666679382324,-7.30000019073486  59.033332824707,-7.03333330154419  59.1666679382324,-6.53333330154419  59.3333320617676,-6.06666660308838  59.6333351135254,-5.43333339691162  59.7333335876465,-4.93333339691162  59.9000015258789,-4.36666679382324  60.0999984741211,-3.79999995231628  60.3333320617676,-3.33333325386047  60.3333320617676,-3.13333344459534  60.4333343505859,-2.43333339691162  60.8333320617676,-0.699999988079071  61.5999984741211,-0.300000011920929  61.6333351135254,0.133333340287209  61.7333335876465,0.699999988079071  61.5999984741211,1.16666662693024  61.7000007629395,3.570012946  62.704987054,5.072916666  63.141654839,5.324345937  63.27982073,5.764583334  63.651352689,5.580264779  63.748901887,5.871012258  64.08125,5.673693034  64.777859701,6.112449529  65.135416667,5.752158435  65.477083333,6.306249999  66.21252281,7.06875  66.601942223,8.008228204  66.81260513,8.927083334  67.329302251,10.239583332  68.009456131,12.09711365  68.623719685,11.8000001907349  68.8000030517578,11.1333332061768  68.6666641235352,9.86666679382324  68.4000015258789,8.86666679382324  67.9000015258789,8.60000038146973  67.9000015258789,8.30000019073486  67.7666702270508,7.76666688919067  67.9000015258789,7.30000019073486  68.2333297729492,6.46666669845581  68.2333297729492,5.76666688919067  68.2333297729492,5.23333311080933  68.4333343505859,4.56666660308838  68.2333297729492,3.83333325386047  68.3666687011719,3.20000004768372  68.3333358764648,1.96666669845581  67.8333358764648,1.515144101  67.800477996,1.5  67.466667175293,1.43333327770233  66.966667175293,1.79999995231628  66.5,2.29999995231628  66.1333312988281,2.06666660308838  65.966667175293,2.90000009536743  65.4333343505859,2.79999995231628  64.9000015258789,3.23333334922791  64.7666702270508,3.63333344459534  64.4333343505859,2.46666669845581  64.1333312988281,1.96666669845581  64.1666641235352,1.014696942  63.626922445,-0.533333361148834  63.4000015258789,-1.73333334922791  63.2666664123535,-2.36791519  63.040542854,-3.20000004768372  63.1333351135254,-4.081391907  62.976460057,-5.76666688919067  63.0666656494141,-5.83333349227905  63.2666664123535,-7.35416173934937  63.6875,-8.39999961853027  63.9333343505859,-10.2666664123535  64.033332824707,-11.4791612625122  64.5208358764648,-12.5333337783813  64.4333343505859,-12.7333335876465  64.2333297729492,-14.2666664123535  63.7333335876465,-15.4333333969116  63.466667175293,-16.2333335876465  63.3333320617676,-17.6333332061768  63.2999992370605,-18.8999996185303  63.2666664123535,-19.7000007629395  63.3666648864746,-20.2333335876465  63.2333335876465,-21.9333324432373  63.3333320617676,-22.466667175293  63.6666679382324,-23.0666675567627  63.2999992370605,-25.7999992370605  60.5))'
$ws4.Range("A94").Value = "WGINOREIW"
$ws4.Range("B94").Value = 'POLYGON((-7.103168115  65.283694796,-6.96666669845581  65.7666702270508,-7.40000009536743  66.0999984741211,-7.330920869  66.384246437,-7.63333320617676  67,-7.519915283  67.25104165,-7.643573575  67.771548904,-7.16666650772095  68.033332824707,-6.93333339691162  68.3666687011719,-7.40000009536743  68.6333312988281,-7.848642619  68.677507828,-7.628274775  68.88475312,-7.561650198  69.247807348,-7.46864723  69.504553404,-7.346337391  69.867374303,-7.167894825  70.118813088,-6.897805367  70.446951911,-6.618760544  70.650056588,-6.173734426  70.746620275,-6.070436289  70.796587087,-6.482534667  70.877686691,-6.968045263  70.957275316,-7.329200449  71.01525725,-7.754596142  71.083122574,-8.283196925  71.153799736,-9.336997827  71.323403555,-10.543476003  71.485203733,-10.770921247  71.510746581,-10.978849338  71.521918122,-11.279133026  71.268840323,-11.787762901  70.977924869,-12.164007593  70.696499491,-12.661247931  70.313189626,-13.025778956  69.965773889,-13.427959537  69.502351037,-13.771387507  69.039790315,-14.145406637  68.687433419,-14.490441666  68.440835791,-14.887272881  68.176801079,-15.199850053  67.925055055,-15.381273911  67.743741963,-15.539809351  67.553362203,-15.591544719  67.381894447,-15.348079034  67.338447774,-14.936077665  67.277714007,-14.46305655  67.21204519,-13.990903778  67.096233947,-13.54589758  66.966882597,-13.268923886  66.853540006,-12.821631959  66.680093502,-12.467434247  66.534149424,-12.267346945  66.42551068,-11.910599311  66.187157192,-11.767199321  66.047561153,-11.611181822  65.837857276,-11.453189531  65.552736547,-12.0333337783813  64.9000015258789,-10.4666662216187  65.0999984741211,-8.89999961853027  65.1666641235352,-7.03333330154419  64.8666687011719,-7.103168115  65.283694796))'
$ws4.Range("A95").Value = "WGINORFaroe"
$ws4.Range("B95").Value = 'POLYGON((-5.93749475479126  62.6875,-5.76666688919067  62.6666679382324,-5.63333320617676  62.5999984741211,-5.43333339691162  62.5666656494141,-5.09999990463257  62.5,-4.83333349227905  62.2999992370605,-4.63333320617676  62.2333335876465,-4.5  62,-4.80000019073486  61.9333343505859,-4.86666679382324  61.9333343505859,-5.03333330154419  61.7999992370605,-4.96666669845581  61.6666679382324,-5.06666660308838  61.6333351135254,-5.23333311080933  61.5666656494141,-5.23333311080933  61.466667175293,-5.30000019073486  61.4333343505859,-5.56666660308838  61.3666648864746,-5.56666660308838  61.2666664123535,-5.80000019073486  61.1666679382324,-5.90000009536743  60.966667175293,-6.09999990463257  60.9000015258789,-6.23333311080933  60.7666664123535,-6.19999980926514  60.6666679382324,-6.36666679382324  60.6333351135254,-6.53333330154419  60.7000007629395,-6.63333320617676  60.7666664123535,-6.83333349227905  60.7999992370605,-7.06666660308838  61.2000007629395,-7.09999990463257  61.2666664123535,-7.23333311080933  61.2999992370605,-7.26666688919067  61.3666648864746,-7.53333330154419  61.4000015258789,-7.56666660308838  61.466667175293,-7.69999980926514  61.466667175293,-7.93333339691162  61.466667175293,-8  61.5,-8  61.533332824707,-7.80000019073486  61.6333351135254,-7.83333349227905  61.7333335876465,-8.10000038146973  61.5999984741211,-8.16666698455811  61.6333351135254,-8.23333358764648  61.7333335876465,-8.33333301544189  61.7999992370605,-8.5  62.2000007629395,-8.43333339691162  62.2999992370605,-8.03333377838135  62.5666656494141,-7.66666650772095  62.6333351135254,-7.43333339691162  62.7000007629395,-7  62.7333335876465,-6.68749475479126  62.7708320617676,-6.53333330154419  62.7333335876465,-6.36666679382324  62.7333335876465,-6.13333320617676  62.7000007629395,-5.93749475479126  62.6875))'
$ws4.Range("A96").Value = "WGINORGS"
$ws4.Range("B96").Value = 'POLYGON((-6.070436289  70.796587087,-6.482534667  70.877686691,-6.968045263  70.957275316,-7.329200449  71.01525725,-7.754596142  71.083122574,-8.283196925  71.153799736,-9.336997827  71.323403555,-10.543476003  71.485203733,-10.770921247  71.510746581,-10.550723259  71.765316178,-10.550723259  71.765316178,-10.361204118  71.992495126,-10.142397394  72.24962139,-9.92533692  72.467780947,-9.741228528  72.652811164,-9.494828604  72.965064819,-9.294028921  73.221165453,-9.155537543  73.368191888,-9.02616899  73.468268821,-8.835775523  73.544080869,-8.500203356  73.612516626,-8.258427496  73.660296347,-7.952277666  73.767075478,-7.574273348  73.894276457,-7.05432894  74.086796099,-6.546865885  74.277219974,-5.970870043  74.50716389,-5.360861181  74.737705597,-4.631983537  74.975532621,-4.144506801  75.143362945,7.752270412  75.340237441,7.975402441  75.06861306,8.184011982  74.910108507,8.40585101  74.747370237,8.631050445  74.585972628,8.773848069  74.443476101,8.827011983  74.298078517,8.76684331  74.104471322,8.623431125  73.927702414,8.383535587  73.773559199,8.019750502  73.563448151,7.788696112  73.44301858,7.423689393  73.338177245,6.766642418  73.22211733,6.121367426  73.062193149,5.596132545  72.97936504,4.883578177  72.891911925,3.748643142  72.69982499,2.945268203  72.578559283,2.001079927  72.415896356,1.46493148  72.323468392,0.883798041  72.2192928,0.141938058  72.069809872,-0.668928933  71.92670097,-1.592431116  71.748990349,-2.738903675  71.520165678,-3.293365895  71.401138856,-3.895005212  71.264388636,-4.350597262  71.162837942,-4.743731131  71.082796373,-5.274472302  70.963746041,-5.758260165  70.863037219,-6.070436289  70.796587087))'
$ws4.Range("A97").Value = "WGINORIFF"
$ws4.Range("B97").Value = 'POLYGON((-7.03333330154419  64.8666687011719,-8.89999961853027  65.1666641235352,-10.4666662216187  65.0999984741211,-12.0333337783813  64.9000015258789,-12.5333337783813  64.4333343505859,-11.4791612625122  64.5208358764648,-10.2666664123535  64.033332824707,-8.39999961853027  63.9333343505859,-7.35416173934937  63.6875,-5.83333349227905  63.2666664123535,-5.76666688919067  63.0666656494141,-4.081391907  62.976460057,-3.20000004768372  63.1333351135254,-4.02082824707031  63.4375,-5.13333320617676  63.9000015258789,-6.26666688919067  64.5666656494141,-6.7  64.8,-7.03333330154419  64.8666687011719))'
$ws4.Range("A98").Value = "WGINORLBE"
$ws4.Range("B98").Value = 'POLYGON((1.5  67.466667175293,1.515144101  67.800477996,1.96666669845581  67.8333358764648,3.20000004768372  68.3333358764648,3.83333325386047  68.3666687011719,4.56666660308838  68.2333297729492,5.23333311080933  68.4333343505859,5.76666688919067  68.2333297729492,6.46666669845581  68.2333297729492,7.30000019073486  68.2333297729492,7.76666688919067  67.9000015258789,8.30000019073486  67.7666702270508,8.60000038146973  67.9000015258789,8.86666679382324  67.9000015258789,9.86666679382324  68.4000015258789,11.1333332061768  68.6666641235352,11.8000001907349  68.8000030517578,12.327593378  68.67342671,12.648402991  68.709930343,12.8333330154419  68.8666687011719,13.371184272  68.925350938,13.731250001  69.082111215,14.25416831  69.122916667,14.506249999  69.190794496,14.861201882  69.255464785,15.514583332  69.440923148,15.9333333969116  69.466667175293,15.994239678  69.672426989,16.637424576  69.747916667,16.63879006  69.940376607,17.113015498  70.11875,17.209717284  70.327083334,17.182235605  70.910416667,16.552674187  71.490174186,16.263299063  71.688299062,15.856249999  71.867016382,15.952789069  72.163877596,16.270093331  72.235416667,16.05316706  72.407333726,15.224489626  72.461989626,15.11196846  72.700531539,15.01421508  73.239951588,15.30480778  73.261858887,15.264583335  73.400843355,15.694392061  73.844392062,16.214601612  74.11875,16.341643079  74.333356919,16.212712404  74.552083333,15.721341619  74.777083333,15.69375  75.006663194,15.368200424  75.134867091,14.871538475  75.165622112,13.973911064  75.199957161,12.030793676  75.261948269,11.968745274  75.106212627,11.915234471  74.749400077,11.826096704  74.454863315,11.644879766  73.872128243,11.362924137  73.218688231,11.165828706  72.855322789,10.943870567  72.459569798,10.794924098  72.248730163,10.626805971  71.984137833,10.491477013  71.830333417,10.33058445  71.699409272,10.006038075  71.514972427,9.837610204  71.444736879,9.531380572  71.325667246,9.242551124  71.221322571,8.837348893  71.110147061,8.527687866  71.008252169,8.202446394  70.934777169,7.875999752  70.864356802,7.462610606  70.776209995,6.956383136  70.669049123,6.365260053  70.552534022,5.963098339  70.472944989,5.434894174  70.369452966,5.218337345  70.321213859,4.864607487  70.266105679,4.455685401  70.218997783,4.166881572  70.181028715,3.406925459  70.081594977,2.712812343  69.995235445,2.467142507  69.947903948,2.207309985  69.900198872,1.739997151  69.829330258,1.09302014  69.706246024,0.364345635  69.562341078,0.032981679  69.474308737,-0.138061852  69.404919172,-0.281662134  69.328212599,-0.40735737  69.223607294,-0.530787038  69.074130028,-0.982562683  68.661066205,-0.696743444  68.503810286,-0.199582561  68.17325349,0.280413613  67.888611434,0.948103416  67.442596719,1.5  67.466667175293))'
$ws4.Range("A99").Value = "WGINORLBW"
$ws4.Range("B99").Value = 'POLYGON((-0.530787038  69.074130028,-0.40735737  69.223607294,-0.281662134  69.328212599,-0.138061852  69.404919172,0.032981679  69.474308737,0.364345635  69.562341078,1.09302014  69.706246024,1.739997151  69.829330258,2.207309985  69.900198872,2.467142507  69.947903948,2.712812343  69.995235445,3.406925459  70.081594977,4.166881572  70.181028715,4.455685401  70.218997783,4.864607487  70.266105679,5.218337345  70.321213859,5.434894174  70.369452966,5.963098339  70.472944989,6.365260053  70.552534022,6.956383136  70.669049123,7.462610606  70.776209995,7.875999752  70.864356802,8.202446394  70.934777169,8.527687866  71.008252169,8.837348893  71.110147061,9.242551124  71.221322571,9.531380572  71.325667246,9.837610204  71.444736879,10.006038075  71.514972427,10.33058445  71.699409272,10.491477013  71.830333417,10.626805971  71.984137833,10.794924098  72.248730163,10.943870567  72.459569798,11.165828706  72.855322789,11.362924137  73.218688231,11.644879766  73.872128243,11.826096704  74.454863315,11.915234471  74.749400077,11.968745274  75.106212627,12.030793676  75.261948269,9.25664  75.3222,7.752270412  75.340237441,7.975402441  75.06861306,8.184011982  74.910108507,8.40585101  74.747370237,8.631050445  74.585972628,8.773848069  74.443476101,8.827011983  74.298078517,8.76684331  74.104471322,8.623431125  73.927702414,8.383535587  73.773559199,8.019750502  73.563448151,7.788696112  73.44301858,7.423689393  73.338177245,6.766642418  73.22211733,6.121367426  73.062193149,5.596132545  72.97936504,4.883578177  72.891911925,3.748643142  72.69982499,2.945268203  72.578559283,2.001079927  72.415896356,1.46493148  72.323468392,0.883798041  72.2192928,0.141938058  72.069809872,-0.668928933  71.92670097,-1.592431116  71.748990349,-2.738903675  71.520165678,-3.293365895  71.401138856,-3.895005212  71.264388636,-4.350597262  71.162837942,-4.743731131  71.082796373,-5.274472302  70.963746041,-5.758260165  70.863037219,-6.070436289  70.796587087,-6.173734426  70.746620275,-5.81072672  70.627767876,-4.954715102  70.327413848,-3.793607086  69.878473848,-2.91985762  69.546454648,-2.386942432  69.340861017,-1.834258608  69.059340681,-1.370058969  68.836437393,-0.982562683  68.661066205,-0.530787038  69.074130028))'
$ws4.Range("A100").Value = "WGINORNorSeaE"
$ws4.Range("B100").Value = 'POLYGON((-0.982562683  68.661066205,-1.372952597  68.274789893,-1.792022466  67.840037576,-1.951116474  67.594830852,-2.154207543  67.384031737,-2.4073482  67.163109131,-2.876927235  66.799914292,-3.482264601  66.298053457,-3.89506653  65.975907178,-4.258846121  65.72917545,-4.736486892  65.52601846,-5.398148656  65.175025901,-5.72851627  65.000144825,-6.7  64.8,-6.26666688919067  64.5666656494141,-5.13333320617676  63.9000015258789,-4.02082824707031  63.4375,-3.20000004768372  63.1333351135254,-2.36791519  63.040542854,-1.73333334922791  63.2666664123535,-0.533333361148834  63.4000015258789,1.014696942  63.626922445,1.96666669845581  64.1666641235352,2.46666669845581  64.1333312988281,3.63333344459534  64.4333343505859,3.23333334922791  64.7666702270508,2.79999995231628  64.9000015258789,2.90000009536743  65.4333343505859,2.06666660308838  65.966667175293,2.29999995231628  66.1333312988281,1.79999995231628  66.5,1.43333327770233  66.966667175293,1.5  67.466667175293,0.948103416  67.442596719,0.280413613  67.888611434,-0.199582561  68.17325349,-0.696743444  68.503810286,-0.982562683  68.661066205))'
$ws4.Range("A101").Value = "WGINORNorSeaW"
$ws4.Range("B101").Value = 'POLYGON((-7.103168115  65.283694796,-6.96666669845581  65.7666702270508,-7.40000009536743  66.0999984741211,-7.330920869  66.384246437,-7.63333320617676  67,-7.519915283  67.25104165,-7.643573575  67.771548904,-7.16666650772095  68.033332824707,-6.93333339691162  68.3666687011719,-7.40000009536743  68.6333312988281,-7.848642619  68.677507828,-7.628274775  68.88475312,-7.561650198  69.247807348,-7.46864723  69.504553404,-7.346337391  69.867374303,-7.167894825  70.118813088,-6.897805367  70.446951911,-6.618760544  70.650056588,-6.173734426  70.746620275,-5.81072672  70.627767876,-4.954715102  70.327413848,-3.793607086  69.878473848,-2.91985762  69.546454648,-2.386942432  69.340861017,-1.834258608  69.059340681,-1.370058969  68.836437393,-0.982562683  68.661066205,-1.372952597  68.274789893,-1.792022466  67.840037576,-1.951116474  67.594830852,-2.154207543  67.384031737,-2.4073482  67.163109131,-2.876927235  66.799914292,-3.482264601  66.298053457,-3.89506653  65.975907178,-4.258846121  65.72917545,-4.736486892  65.52601846,-5.398148656  65.175025901,-5.72851627  65.000144825,-6.7  64.8,-7.03333330154419  64.8666687011719,-7.103168115  65.283694796))'
$ws4.Range("A102").Value = "WGINORRH"
$ws4.Range("B102").Value = 'POLYGON((-14.0666666030884  58.033332824707,-14.5333337783813  57.6666679382324,-14.7666664123535  57.1666679382324,-15.1666669845581  56.7000007629395,-15.6666669845581  56.2333335876465,-15.5333337783813  56.033332824707,-14.9666662216187  56.0666656494141,-14.2666664123535  56.5666656494141,-13.4333333969116  57,-12.9333333969116  57.4000015258789,-13.1333332061768  57.9000015258789,-13.6333332061768  58.1666679382324,-14.0666666030884  58.033332824707))'

# --- Selections / active view state to match target ---
$ws1.Activate()
$ws1.Range("G8").Select()
$ws2.Activate()
$ws2.Range("D94").Select()
$ws4.Activate()
$ws4.Range("B102").Select()
$ws1.Activate()
